# Add a new worksheet "2020-11-10" at the end of the workbook (after the
# current last sheet) and populate it with the attendance/heart-rate data,
# matching the "updated the heartrate into the gui" commit.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "2020-11-10"

# Header row
$headers = @("Sr. No", "Name", "Address", "Job", "Time-Stamp", "SpO2_value")
for ($col = 1; $col -le $headers.Length; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headers[$col - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Data rows: Time-Stamp, SpO2_value, (optional) heart-rate value in column G
$data = @(
    @("01:38:44", 84.61483623004834, $null),
    @("01:41:23", 82.09468127949228, $null),
    @("09:21:53", 89.30021772238722, $null),
    @("09:35:17", 94.56050200781256, $null),
    @("09:36:59", 92.57248194477947, $null),
    @("09:43:12", 94.81811911912854, $null),
    @("09:47:16", 96.61184229430469, $null),
    @("09:53:50", 94.65548776817448, $null),
    @("09:54:29", 95.39734931207116, $null),
    @("09:59:30", 95.39645256607388, 66.23650871842609),
    @("12:10:10", 95.20561567862042, 58.05271621944848),
    @("12:14:45", 96.11290914378348, 74.2404288587264),
    @("12:17:37", 90.97632885109677, 75.00742923083941),
    @("12:27:34", 92.9750562342219, 65.77138664733151),
    @("12:30:34", 97.15916212822721, 54.99891500210524)
)

$row = 2
foreach ($item in $data) {
    $aCell = $ws.Cells.Item($row, 1)
    $aCell.Value = 1
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.LineStyle = 1

    $ws.Cells.Item($row, 2).Value = "sachin"
    $ws.Cells.Item($row, 3).Value = "xyz/xyz"
    $ws.Cells.Item($row, 4).Value = "coder"
    $ws.Cells.Item($row, 5).Value = $item[0]
    $ws.Cells.Item($row, 6).Value = $item[1]
    if ($item[2] -ne $null) {
        $ws.Cells.Item($row, 7).Value = $item[2]
    }

    $row = $row + 1
}

[void]$ws.Range("A1").Select()
[void]$lastSheet.Activate()

Write-Host "Added sheet 2020-11-10 with $($data.Length) data rows"
